{"js": "// Each entry's `oldText` is the full, original text of one \"discussion\"\n// paragraph's trailing run; `newText` is that same paragraph after the\n// author inserted extra sentences (per the commit's diff).\nconst replacements = [\n  {\n    oldText:\n      \"The Glioblastoma Cancer shows abberations in copy number in Chromosome 7, Chromosome 19 and Chromosome 20 as a gain and in copy 9p, Chromosome 10, chromosome 13q, Chromosome 14q, Chromosome 21q as a loss. Chromosome 13, 14 and 15 show very little abberation in the p arm.\",\n    newText:\n      \"The Glioblastoma Cancer type is often the most agressive type of cancer that begins within the brain. It shows abberations in copy number in Chromosome 7, Chromosome 19 and Chromosome 20 as a gain and in copy 9p, Chromosome 10, chromosome 13q, Chromosome 14q, Chromosome 21q as a loss. Chromosome 13, 14 and 15 show very little abberation in the p arm. The literature conferms alteration in Chromosome 7, 9 and 10 to be the most significant.\",\n  },\n  {\n    oldText:\n      \"The invasive breast cancer carsinoma shows CNVs in Chromosome 1 in the q arm, in Chromosome 8 in the q arm a copy number gain. In Chromosome 1 in the p arm, Chromosome 4, Chromosome 16q and Chromosome X a loss.\",\n    newText:\n      \"The invasive breast cancer carsinoma shows CNVs in Chromosome 1 in the q arm, in Chromosome 8 in the q arm a copy number gain. Chromosome 5 could be interesting. chromosome In Chromosome 1 in the p arm, Chromosome 4, Chromosome 16q and Chromosome X a loss. Chromosome 16p arm shows gain. In invasive breast cancer, different types show different patterns of chromosomal abberations. The literature supports that Chromosome 16 and X could play a vital role in breast cancer.\",\n  },\n  {\n    oldText:\n      \"The Lung non-small call carcinoma shows a lot of abberation in Chromosome 1, Chromosome 3p, Chromosome 5p, Chromosome 8p and q (loss and gain), and apperantly a complete loss of Chromosome 13p, 14p, 15p, 21p, 22p. Again here, Chromosome 13, 14 and 15 show again very little change in the p arm.\",\n    newText:\n      \"The Lung non-small call carcinoma shows a lot of abberation in Chromosome 1, Chromosome 3p, Chromosome 5p, Chromosome 8p and q (loss and gain), and apperantly a complete loss of Chromosome 13p, 14p, 15p, 21p, 22p. Again here, Chromosome 13, 14 and 15 show again very little change in the p arm. It is supported in the literature that in NSCC chromosome 5, 7 and 8 show a gain of copy numbers.\",\n  },\n  {\n    oldText:\n      \"The colon adenocarcinoma shows abberations in Chromosome 4, 7, 8, 13, 18, 20. For Chromosome 14p, it could be that the data gathered is not enough or there are no probes detecting the 14p arm. For the q arm however, there seems to be a trend for a loss of an allel. Interesting is that in general the background noise in this Cancer is lower and the extremer abberations are more prominent.\",\n    newText:\n      \"The colon adenocarcinoma is one of the most common inherited cancer syndromes known.shows abberations in Chromosome 4, 7, 8, 13, 18, 20. For Chromosome 14p, it could be that the data gathered is not enough or there are no probes detecting the 14p arm. For the q arm however, there seems to be a trend for a loss of an allel. Interesting is that in general the background noise in this Cancer is lower and the extremer abberations are more prominent. In the literature Chromosome 1p, 5q, 8p, 15q, .18 q (which is very significant in this data) have been set in context with colon cancer.\",\n  },\n  {\n    oldText:\n      \"In the Melanoma especially Chromosomes 4, 7, 9 10, maybe 11, 13, 14, 15, 21 and 22 show trends of abnormal. Samples 13, 14, 15 all show very low aberation in the p arm, so the noise is cancelled out. The changes are more fluctuating, there are single positions within, which show stronger abberations.\",\n    newText:\n      \"In the Melanoma especially Chromosomes 4, 6q, 7, 9 10, maybe 11, 13, 14, 15, 21 and 22 show trends of abnormal. Samples 13, 14, 15 all show very low aberation in the p arm, so the noise is cancelled out. The changes are more fluctuating, there are single positions within, which show stronger abberations. In the literature it has been described that Chromosome 6, 7, 9, 10 play a role in Melanoma cancer. Apparentlyz abberations in 9 and 10 occur early on in melanoma progression, whereas gains of chromosome 7 occur later.\",\n  },\n];\n\nconst body = context.document.body;\nlet count = 0;\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText.slice(0, 40));\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  count++;\n}\nawait context.sync();\nreturn \"replaced=\" + count;", "ps1": "# Expands five \"discussion\" paragraphs (Glioblastoma, Invasive Breast\n# Carcinoma, Lung NSCC, Colon Adenocarcinoma, Melanoma) with extra\n# sentences the author inserted, per the commit's diff. Each paragraph's\n# entire original text is matched exactly once and replaced in place.\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($doc, $findText, $replaceText) {\n  $rng = $doc.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  # wdFindContinue=1 (restrict to one exact match), wdReplaceOne=1\n  $found = $rng.Find.Execute($findText, $false, $true, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n  if (-not $found) {\n    throw \"Text not found: $findText\"\n  }\n}\n\n$find = \"The Glioblastoma Cancer shows abberations in copy number in Chromosome 7, Chromosome 19 and Chromosome 20 as a gain and in copy 9p, Chromosome 10, chromosome 13q, Chromosome 14q, Chromosome 21q as a loss. Chromosome 13, 14 and 15 show very little abberation in the p arm.\"\n$replace = \"The Glioblastoma Cancer type is often the most agressive type of cancer that begins within the brain. It shows abberations in copy number in Chromosome 7, Chromosome 19 and Chromosome 20 as a gain and in copy 9p, Chromosome 10, chromosome 13q, Chromosome 14q, Chromosome 21q as a loss. Chromosome 13, 14 and 15 show very little abberation in the p arm. The literature conferms alteration in Chromosome 7, 9 and 10 to be the most significant.\"\nReplace-ExactText $d $find $replace\n\n$find = \"The invasive breast cancer carsinoma shows CNVs in Chromosome 1 in the q arm, in Chromosome 8 in the q arm a copy number gain. In Chromosome 1 in the p arm, Chromosome 4, Chromosome 16q and Chromosome X a loss.\"\n$replace = \"The invasive breast cancer carsinoma shows CNVs in Chromosome 1 in the q arm, in Chromosome 8 in the q arm a copy number gain. Chromosome 5 could be interesting. chromosome In Chromosome 1 in the p arm, Chromosome 4, Chromosome 16q and Chromosome X a loss. Chromosome 16p arm shows gain. In invasive breast cancer, different types show different patterns of chromosomal abberations. The literature supports that Chromosome 16 and X could play a vital role in breast cancer.\"\nReplace-ExactText $d $find $replace\n\n$find = \"The Lung non-small call carcinoma shows a lot of abberation in Chromosome 1, Chromosome 3p, Chromosome 5p, Chromosome 8p and q (loss and gain), and apperantly a complete loss of Chromosome 13p, 14p, 15p, 21p, 22p. Again here, Chromosome 13, 14 and 15 show again very little change in the p arm.\"\n$replace = \"The Lung non-small call carcinoma shows a lot of abberation in Chromosome 1, Chromosome 3p, Chromosome 5p, Chromosome 8p and q (loss and gain), and apperantly a complete loss of Chromosome 13p, 14p, 15p, 21p, 22p. Again here, Chromosome 13, 14 and 15 show again very little change in the p arm. It is supported in the literature that in NSCC chromosome 5, 7 and 8 show a gain of copy numbers.\"\nReplace-ExactText $d $find $replace\n\n$find = \"The colon adenocarcinoma shows abberations in Chromosome 4, 7, 8, 13, 18, 20. For Chromosome 14p, it could be that the data gathered is not enough or there are no probes detecting the 14p arm. For the q arm however, there seems to be a trend for a loss of an allel. Interesting is that in general the background noise in this Cancer is lower and the extremer abberations are more prominent.\"\n$replace = \"The colon adenocarcinoma is one of the most common inherited cancer syndromes known.shows abberations in Chromosome 4, 7, 8, 13, 18, 20. For Chromosome 14p, it could be that the data gathered is not enough or there are no probes detecting the 14p arm. For the q arm however, there seems to be a trend for a loss of an allel. Interesting is that in general the background noise in this Cancer is lower and the extremer abberations are more prominent. In the literature Chromosome 1p, 5q, 8p, 15q, .18 q (which is very significant in this data) have been set in context with colon cancer.\"\nReplace-ExactText $d $find $replace\n\n$find = \"In the Melanoma especially Chromosomes 4, 7, 9 10, maybe 11, 13, 14, 15, 21 and 22 show trends of abnormal. Samples 13, 14, 15 all show very low aberation in the p arm, so the noise is cancelled out. The changes are more fluctuating, there are single positions within, which show stronger abberations.\"\n$replace = \"In the Melanoma especially Chromosomes 4, 6q, 7, 9 10, maybe 11, 13, 14, 15, 21 and 22 show trends of abnormal. Samples 13, 14, 15 all show very low aberation in the p arm, so the noise is cancelled out. The changes are more fluctuating, there are single positions within, which show stronger abberations. In the literature it has been described that Chromosome 6, 7, 9, 10 play a role in Melanoma cancer. Apparentlyz abberations in 9 and 10 occur early on in melanoma progression, whereas gains of chromosome 7 occur later.\"\nReplace-ExactText $d $find $replace\n\nWrite-Output \"done\""}
